$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Customer success" / "inbox" column (column C) entirely -
#    clearing its contents drops it from sheetData/dimension/shared strings.
$ws.Range("C1:C2").ClearContents()

# 2. Update the credential values shown in the sheet.
$ws.Range("A2").Value = "Test_UX09@westpharma.com"
$ws.Range("B2").Value = "Westpharm@2019"

# 3. Add a hyperlink on the password cell (B2), matching the pattern
#    already used for the email cell (A2).
$origStyle = $ws.Range("A2").Style()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Test_UX09@westpharma.com")

# Re-apply the original (shared) Hyperlink cell style to B2 so it reuses
# the same style index as A2 instead of a freshly minted one.
$ws.Range("B2").Style = $origStyle

# 4. Widen column B slightly to fit the new, longer password value.
$ws.Columns("B").ColumnWidth = 15.67
